# Excel_models.xlsx edit script
# 1. Rename the worksheet tab from "Feuil1" to "Excel_models"
# 2. Add a new "tester" column (D) with reviewer names assigned to each model row
# 3. Reset the view selection to E1 (also clears the stale topLeftCell scroll state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Excel_models"

# Tester assignment per row, cycling farnaz / matevz / dorian / ali / gio / raiko
$testers = @{
    2  = "farnaz"
    3  = "matevz"
    4  = "dorian"
    5  = "ali"
    6  = "gio"
    7  = "raiko"
    8  = "farnaz"
    9  = "matevz"
    10 = "dorian"
    11 = "ali"
    12 = "gio"
    13 = "raiko"
    14 = "farnaz"
    15 = "matevz"
    16 = "dorian"
    17 = "ali"
    18 = "gio"
    19 = "raiko"
    20 = "dorian"
    21 = "dorian"
    22 = "dorian"
    23 = "dorian"
}

# Rows that don't yet have a D cell formatted like the rest of column D (s="2")
# need their formatting copied from an existing formatted D cell (D2) before
# the value is written.
$needsFormat = @(10, 14, 15, 16, 17, 18, 19, 20, 21, 22, 23)

$ws.Range("D2").Copy() | Out-Null
foreach ($row in $needsFormat) {
    $ws.Cells.Item($row, 4).PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = $false

# Header for the new column (no special formatting, like the rest of row 1)
$ws.Range("D1").Value = "tester"

foreach ($row in $testers.Keys) {
    $ws.Cells.Item($row, 4).Value = $testers[$row]
}

# Update the view: select E1 (this also clears the old topLeftCell scroll anchor)
$ws.Range("E1").Select() | Out-Null
